# Hassu "migraatiotesti" workbook edit:
#  - Add a new sheet "Migraatio2" between "Migraatio" and "Metadata", which
#    takes over the migrated data row that used to live in row 5 of "Migraatio".
#  - Clear that row's content (A5:D5) on the original "Migraatio" sheet, leaving
#    only the still-styled, empty D5 cell behind.
#  - Leave "Metadata" untouched (it just shifts from 2nd to 3rd tab).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Migraatio")

# Duplicate "Migraatio" (keeps column widths / number formats / the x14 data
# validation / everything) and drop it immediately after the original sheet,
# then rename the duplicate to "Migraatio2".
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Migraatio2"

# "Migraatio2" should only contain the header row plus the single migrated
# record (originally on row 5) and a trailing blank row, so remove the rows
# that held the other sample records.
$ws2.Rows("2:4").Delete()

# Remove the migrated record from the original "Migraatio" sheet - row 5
# stays in place (dimension/rows below are untouched), just emptied.
$ws1.Range("A5:D5").ClearContents()

# Match the saved selections / active tab from the edit: row 5 is selected
# (as if it had just been cut) on "Migraatio", "Migraatio2" has cell A13
# selected and is the active/visible sheet.
$ws1.Rows("5:5").Select()
$ws2.Range("A13").Select()
$ws2.Activate()
